# Updates cryptos list (price / 1h volume columns) and refreshes a few
# coin rows (47-51) to reflect the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.641.68"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "2.487.13"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'585.91"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "'176.31"
$ws.Range("E6").Value = "  +4.87%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("E9").Value = "  +4.12%  "
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("E11").Value = "  +3.07%  "
$ws.Range("D12").Value = "'4.93"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "2.949.86"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "67.539.14"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").Value = "2.462.03"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "'11.05"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "'7.40"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").Value = "'351.62"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "'4.04"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'70.56"
$ws.Range("E23").Value = "  +2.78%  "
$ws.Range("D24").Value = "'4.23"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").Value = "'1.80"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").Value = "'9.22"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").Value = "2.615.29"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").Value = "'0.996"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").Value = "'511.50"
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("D31").Value = "'7.83"
$ws.Range("E31").Value = "  +2.31%  "
$ws.Range("D32").Value = "'1.25"
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  +6.61%  "
$ws.Range("D36").Value = "'160.55"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "'18.34"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("D42").Value = "'0.330"
$ws.Range("E42").Value = "  +1.58%  "
$ws.Range("D43").Value = "'4.87"
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("D44").Value = "'2.43"
$ws.Range("E44").Value = "  +3.37%  "
$ws.Range("D45").Value = "'143.60"
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("D46").Value = "'3.51"
$ws.Range("E46").Value = "  +2.50%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'0.514"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0747"
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("B49").Value = "Optimism"
$ws.Range("C49").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D49").Value = "'1.58"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.586"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D51").Value = "'1.19"
$ws.Range("E51").Value = "  +1.80%  "
